$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Footer 1 (word/footer1.xml): Pearson logo, docPr id="1" ---
$f1 = $sec.Footers(1)
$shp1 = $f1.Range.InlineShapes(1)
$shp1.Name = "image1.png"

# --- Footer 2 (word/footer2.xml): Pearson logo, docPr id="2" ---
$f2 = $sec.Footers(2)
$shp2 = $f2.Range.InlineShapes(1)
$shp2.Name = "image1.png"

# --- Header 2 (word/header2.xml): BTec logo, docPr id="3" ---
$h2 = $sec.Headers(2)
$shp3 = $h2.Range.InlineShapes(1)
$shp3.Name = "image2.jpg"
